$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8644635677337646
$ws.Range("B1").Value = 3.259998798370361
$ws.Range("C1").Value = 2.597519159317017
$ws.Range("D1").Value = 2.43107533454895
$ws.Range("E1").Value = 2.026140928268433
